$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.065.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.502.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.38"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.44"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.54"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.20"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +7.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0810"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.08"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.893.03"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.499.69"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "47.934.34"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.94"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.71"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0938"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.49%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.31"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +12.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.47"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.53"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.84"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.66"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.41%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.30"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.09"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.49"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.48"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0782"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.62"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "120.87"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.16"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.73%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.017.40"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.99"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.84"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.05"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.96%  "
